$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.460.19'
$ws.Range('E2').Value = '  +0.39%  '
$ws.Range('D3').Value = '1.569.55'
$ws.Range('E3').Value = '  -1.50%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').Value = '211.81'
$ws.Range('E5').Value = '  -1.22%  '
$ws.Range('D6').Value = '0.492'
$ws.Range('E6').Value = '  -0.46%  '
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('D8').Value = '45.99'
$ws.Range('E8').Value = '  +4.70%  '
$ws.Range('D9').Value = '23.97'
$ws.Range('E9').Value = '  -0.37%  '
$ws.Range('E10').Value = '  -1.84%  '
$ws.Range('D11').Value = '0.0589'
$ws.Range('E11').Value = '  -1.84%  '
$ws.Range('E12').Value = '  +0.03%  '
$ws.Range('D13').Value = '1.793.51'
$ws.Range('E13').Value = '  -1.61%  '
$ws.Range('D14').Value = '1.564.54'
$ws.Range('E14').Value = '  -1.81%  '
$ws.Range('D15').Value = '0.519'
$ws.Range('E15').Value = '  -2.28%  '
$ws.Range('D16').Value = '28.427.97'
$ws.Range('E16').Value = '  +0.19%  '
$ws.Range('E17').Value = '  -2.36%  '
$ws.Range('D18').Value = '62.07'
$ws.Range('E18').Value = '  -1.66%  '
$ws.Range('D19').Value = '228.93'
$ws.Range('E19').Value = '  +0.75%  '
$ws.Range('E20').Value = '  -2.20%  '
$ws.Range('D21').Value = '0.0₃0688'
$ws.Range('E21').Value = '  -3.06%  '
$ws.Range('E22').Value = '  +0.00%  '
$ws.Range('D23').Value = '3.87'
$ws.Range('E23').Value = '  -5.41%  '
$ws.Range('D24').Value = '9.07'
$ws.Range('E24').Value = '  -2.59%  '
$ws.Range('D25').Value = '2.11'
$ws.Range('E25').Value = '  +7.09%  '
$ws.Range('D26').Value = '150.55'
$ws.Range('E26').Value = '  -0.78%  '
$ws.Range('D27').Value = '14.96'
$ws.Range('E27').Value = '  -1.57%  '
$ws.Range('D28').Value = '6.42'
$ws.Range('E28').Value = '  -2.44%  '
$ws.Range('E29').Value = '  -3.76%  '
$ws.Range('E30').Value = '  -0.09%  '
$ws.Range('D31').Value = '0.0477'
$ws.Range('E31').Value = '  +0.50%  '
$ws.Range('D32').Value = '1.11'
$ws.Range('E32').Value = '  -3.52%  '
$ws.Range('E33').Value = '  -1.49%  '
$ws.Range('D34').Value = '3.07'
$ws.Range('E34').Value = '  -1.86%  '
$ws.Range('D35').Value = '1.391.28'
$ws.Range('E35').Value = '  -0.22%  '
$ws.Range('D36').Value = '1.06'
$ws.Range('E36').Value = '  +1.30%  '
$ws.Range('E37').Value = '  -3.40%  '
$ws.Range('E38').Value = '  +0.33%  '
$ws.Range('D39').Value = '2.62'
$ws.Range('E39').Value = '  +4.36%  '
$ws.Range('D40').Value = '0.0165'
$ws.Range('E40').Value = '  -0.98%  '
$ws.Range('D41').Value = '0.523'
$ws.Range('E41').Value = '  -2.97%  '
$ws.Range('E42').Value = '  +0.02%  '
$ws.Range('D43').Value = '1.88'
$ws.Range('E43').Value = '  +0.69%  '
$ws.Range('D44').Value = '0.786'
$ws.Range('E44').Value = '  -3.38%  '
$ws.Range('B45').Value = 'Kaspa'
$ws.Range('C45').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D45').Value = '0.0468'
$ws.Range('E45').Value = '  +4.23%  '
$ws.Range('E46').Value = '  -4.54%  '
$ws.Range('B47').Value = 'WEMIXToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D47').Value = '0.969'
$ws.Range('E47').Value = '  -1.59%  '
$ws.Range('B48').Value = 'Aave'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D48').Value = '62.67'
$ws.Range('E48').Value = '  -2.71%  '
$ws.Range('B49').Value = 'RocketPoolETH'
$ws.Range('C49').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D49').Value = '1.705.69'
$ws.Range('E49').Value = '  -1.60%  '
$ws.Range('B50').Value = 'Quant'
$ws.Range('C50').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D50').Value = '85.93'
$ws.Range('E50').Value = '  -1.73%  '
$ws.Range('B51').Value = 'BabyDogeCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('E51').Value = '  -1.35%  '
